$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 275; existing rows 275-312 shift down to 276-313.
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row with its data.
$ws.Cells.Item(275, 1).Value = 3
$ws.Cells.Item(275, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(275, 3).Value = "Coquimbo"
$ws.Cells.Item(275, 4).Value = [DateTime]"2022-07-06"
$ws.Cells.Item(275, 5).Value = 5
$ws.Cells.Item(275, 6).Value = 100112001
$ws.Cells.Item(275, 7).Value = "Berenjena"
$ws.Cells.Item(275, 8).Value = "Sin especificar"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 105
$ws.Cells.Item(275, 11).Value = 7000
$ws.Cells.Item(275, 12).Value = 8000
$ws.Cells.Item(275, 13).Value = 7524
$ws.Cells.Item(275, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(275, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(275, 16).Value = 125
$ws.Cells.Item(275, 17).Value = 60
$ws.Cells.Item(275, 18).Value = "Hortaliza"
